# "added 4wk low sales check" - update forecast figures on "Forecast Comparison"
# sheet and roll the new totals up into the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---------------------------------------------
# Row 2 (W10)
$wsForecast.Range("D2").Value = 9
$wsForecast.Range("H2").Value = 12.93
$wsForecast.Range("L2").Value = 0.85

# Row 3 (W11)
$wsForecast.Range("D3").Value = 9
$wsForecast.Range("H3").Value = 11.68
$wsForecast.Range("L3").Value = 1.04

# Row 4 (W12)
$wsForecast.Range("D4").Value = 9
$wsForecast.Range("H4").Value = 10.24
$wsForecast.Range("L4").Value = 1.13

# Row 5 (W13)
$wsForecast.Range("D5").Value = 10
$wsForecast.Range("H5").Value = 9.06
$wsForecast.Range("L5").Value = 1.1

# Row 6 (W14)
$wsForecast.Range("H6").Value = 8.06
$wsForecast.Range("L6").Value = 0.85

# Row 7 (W15)
$wsForecast.Range("D7").Value = 10
$wsForecast.Range("H7").Value = 7.06
$wsForecast.Range("L7").Value = 0.81

# Row 8 (W16)
$wsForecast.Range("D8").Value = 10
$wsForecast.Range("H8").Value = 6.06
$wsForecast.Range("L8").Value = 0.89

# Row 9 (W17)
$wsForecast.Range("D9").Value = 10
$wsForecast.Range("H9").Value = 4.96
$wsForecast.Range("L9").Value = 1.11

# Row 10 (W18)
$wsForecast.Range("D10").Value = 10
$wsForecast.Range("H10").Value = 3.88
$wsForecast.Range("L10").Value = 1.07

# Row 11 (W19)
$wsForecast.Range("D11").Value = 10
$wsForecast.Range("H11").Value = 2.88
$wsForecast.Range("L11").Value = 1.19

# Row 12 (W20)
$wsForecast.Range("D12").Value = 10
$wsForecast.Range("H12").Value = 1.85
$wsForecast.Range("L12").Value = 1.14

# Row 13 (W21)
$wsForecast.Range("D13").Value = 10
$wsForecast.Range("H13").Value = 0.85
$wsForecast.Range("J13").Value = "Urgent"
$wsForecast.Range("L13").Value = 1.03

# Row 14 (W22)
$wsForecast.Range("D14").Value = 10
$wsForecast.Range("H14").Value = 0
$wsForecast.Range("I14").Value = "High"
$wsForecast.Range("J14").Value = "Urgent"
$wsForecast.Range("L14").Value = 0.85

# Row 15 (W23)
$wsForecast.Range("D15").Value = 11
$wsForecast.Range("H15").Value = 0
$wsForecast.Range("I15").Value = "High"
$wsForecast.Range("J15").Value = "Urgent"
$wsForecast.Range("L15").Value = 1.08

# Row 16 (W24)
$wsForecast.Range("D16").Value = 11
$wsForecast.Range("H16").Value = 0
$wsForecast.Range("L16").Value = 1.03

# Row 17 (W25)
$wsForecast.Range("D17").Value = 11
$wsForecast.Range("L17").Value = 1.15

# --- Summary sheet ----------------------------------------------------------
# These values are stored as text in the workbook, so force a text number
# format before writing so Excel does not auto-convert them to numbers.
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "164"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "79"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "38"

$wsSummary.Range("B12").NumberFormat = "@"
$wsSummary.Range("B12").Value = "11"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "9"
